# Weekly data refresh: insert two new rows of "Piña" price data at the top
# of the Vega Modelo de Temuco / Caramelo block (rows 259-260), pushing the
# existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 259; this shifts rows 259:313 down
# to 261:315 and carries the row-259 formatting (incl. the date style on
# column D) onto the freshly inserted rows.
$ws.Rows("259:260").Insert()

# --- Row 259: new "Primera" quality entry ---
$ws.Range("A259").Value = 10
$ws.Range("B259").Value = "Vega Modelo de Temuco"
$ws.Range("C259").Value = "La Araucanía"
$ws.Range("D259").Value = 44511
$ws.Range("E259").Value = 9
$ws.Range("F259").Value = "Fruta"
$ws.Range("G259").Value = 100108
$ws.Range("H259").Value = "Tropicales y subtropicales"
$ws.Range("I259").Value = 100108005
$ws.Range("J259").Value = "Piña"
$ws.Range("K259").Value = "Caramelo"
$ws.Range("L259").Value = "Primera"
$ws.Range("M259").Value = 93
$ws.Range("N259").Value = 20000
$ws.Range("O259").Value = 21000
$ws.Range("P259").Value = 20376
$ws.Range("Q259").Value = "$/caja 12 unidades"
$ws.Range("R259").Value = "Ecuador"
$ws.Range("S259").Value = 1698
$ws.Range("T259").Value = 12

# --- Row 260: new "Segunda" quality entry ---
$ws.Range("A260").Value = 10
$ws.Range("B260").Value = "Vega Modelo de Temuco"
$ws.Range("C260").Value = "La Araucanía"
$ws.Range("D260").Value = 44511
$ws.Range("E260").Value = 9
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100108
$ws.Range("H260").Value = "Tropicales y subtropicales"
$ws.Range("I260").Value = 100108005
$ws.Range("J260").Value = "Piña"
$ws.Range("K260").Value = "Caramelo"
$ws.Range("L260").Value = "Segunda"
$ws.Range("M260").Value = 75
$ws.Range("N260").Value = 20000
$ws.Range("O260").Value = 20000
$ws.Range("P260").Value = 20000
$ws.Range("Q260").Value = "$/caja 14 unidades"
$ws.Range("R260").Value = "Ecuador"
$ws.Range("S260").Value = 1429
$ws.Range("T260").Value = 14
